# Generate Report for Handoff
# Adds a new handoff row (row 3) to the Overview, zh-cn, and de-de sheets,
# mirroring the existing row 2 pattern for the new file
# 8f66ff09-ca8c-4f1b-b867-15a1956dbcd5ooo....md

$wb = $excel.ActiveWorkbook

$mdName     = '8f66ff09-ca8c-4f1b-b867-15a1956dbcd5oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdDisplay  = 'e2e\8f66ff09-ca8c-4f1b-b867-15a1956dbcd5oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$ready      = 'Ready for handoff'
$dt1019     = '2017-02-17 10:19:01'
$zhcnXlf    = '8f66ff09-ca8c-4f1b-b867-15a1956dbcd5ooooooooooooooooooooooooooooooooooooo.fb47c82f2a8408d2fb0cb029139eebbd8a706cd3.zh-cn.xlf'
$dedeXlf    = '8f66ff09-ca8c-4f1b-b867-15a1956dbcd5ooooooooooooooooooooooooooooooooooooo.fb47c82f2a8408d2fb0cb029139eebbd8a706cd3.de-de.xlf'
$dt101845   = '2017-02-17 10:18:45'
$newUrl     = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfd1f140332a529d01a96b7ca8e8924f2d517931/e2e/8f66ff09-ca8c-4f1b-b867-15a1956dbcd5oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$dateFmt    = 'yyyy-mm-dd HH:mm:ss'

# ---------------------------------------------------------------------
# Sheet "Overview": add row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", $mdDisplay) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $ready
$wsOverview.Range("F3").Value = $ready
$wsOverview.Range("G3").Value = $dt1019
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))
$wsOverview.Range("E5").ColumnWidth = 17.2159881591797
$wsOverview.Range("F5").ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# Sheet "zh-cn": add row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newUrl, "", "", $mdName) | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $ready
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $dt101845
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = ""
$wsZhCn.Range("L3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").NumberFormat = $dateFmt
$wsZhCn.Range("M3").Value = ""
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "True"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Range("Q3").Value = "False"
$wsZhCn.Range("R3").Value = ""

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:R3"))
$wsZhCn.Range("C5").ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# Sheet "de-de": add row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newUrl, "", "", $mdName) | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $ready
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $dt1019
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = ""
$wsDeDe.Range("L3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").NumberFormat = $dateFmt
$wsDeDe.Range("M3").Value = ""
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "True"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Range("Q3").Value = "False"
$wsDeDe.Range("R3").Value = ""

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:R3"))
$wsDeDe.Range("C5").ColumnWidth = 17.2159881591797

Write-Host "Done: added handoff row for 8f66ff09 across Overview/zh-cn/de-de"
